$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.447.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7139"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07918"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3124"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08284"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.886.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7293"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.301"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.455.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.947"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.05%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "248.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.09%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007878"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.131.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.60%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.980"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.97%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1607"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.039"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.63%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.53%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.364"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.63%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.20%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.412"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.135"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05320"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.938"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.198"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.74%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7272"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.680"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.29%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.242.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.95%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.724"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9117"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.43%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.04%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.192"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.030.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5313"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.970"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.88%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.768"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.321"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.87%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4329"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.56%  "
